# HillTestData.xlsx -- "Making the pandas more smart"
# Give the two data columns descriptive headers (instead of the generic
# "x" / "y") so the pandas/curve-fitting notebook that reads this sheet
# picks up meaningful column names, and tidy up the column widths/
# selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the ligand concentration (lasAHL), column B holds the
# reporter readout (GFP geometric mean MEFL).
$ws.Range("A1").Value = "lasAHL (nM)"
$ws.Range("B1").Value = "GFP (geometric mean MEFL)"

# Widen the header columns so the new, longer labels are fully visible.
$ws.Columns("A:B").AutoFit()

# Leave the selection on the first data row, under the headers.
$ws.Range("A2").Select()
